$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("LP1912")
$ws215  = $wb.Worksheets.Item("LP1912-215")
$ws6203 = $wb.Worksheets.Item("6203-6173")

function Set-DataRow {
    param($Row, $A, $B, $C, $D, $E)
    $wsMain.Cells.Item($Row, 1).Value = $A
    $wsMain.Cells.Item($Row, 2).Value = $B
    $wsMain.Cells.Item($Row, 3).Value = $C
    $wsMain.Cells.Item($Row, 4).Value = $D
    $wsMain.Cells.Item($Row, 5).Value = $E
}

# --- Header rows: "Ultima actualizacion" / "Total filas" ---
$wsMain.Cells.Item(2, 1).Value = "Última actualización: 16:53:02"
$wsMain.Cells.Item(3, 1).Value = "Total filas: 342"
$ws215.Cells.Item(2, 1).Value = "Última actualización: 16:53:02"
$ws6203.Cells.Item(2, 1).Value = "Última actualización: 16:53:02"

# --- Data rows (added / reordered / value-swapped) ---
Set-DataRow 38 '06:44:40' '07:05' '23_HERNANDEZ' 21 'LP1912'
Set-DataRow 39 '05:18:42' '07:05' '15_ABASTO' 107 'LP1912'
Set-DataRow 48 '07:24:45' '07:31' '16_SANTA ANA' 7 'LP1912'
Set-DataRow 49 '05:53:46' '07:31' '11_ETCHEVERRY' 98 'LP1912'
Set-DataRow 90 '08:55:25' '09:16' '16_SANTA ANA' 21 'LP1912'
Set-DataRow 91 '07:24:45' '09:16' '27_EL RETIRO' 112 'LP1912'
Set-DataRow 133 '09:26:30' '11:06' '16_P MOR-167 Y 521' 100 'LP1912'
Set-DataRow 134 '10:52:37' '11:06' '23_HERNANDEZ' 14 'LP1912'
Set-DataRow 165 '12:01:11' '12:07' '16_P MOR-SANTA ANA' 6 'LP1912'
Set-DataRow 166 '11:46:46' '12:07' '23_HERNANDEZ' 21 'LP1912'
Set-DataRow 186 '11:17:39' '12:41' '10_OLMOS' 84 'LP1912'
Set-DataRow 187 '11:46:46' '12:41' '23_HERNANDEZ' 55 'LP1912'
Set-DataRow 208 '11:46:46' '13:26' '14_ABASTO' 100 'LP1912'
Set-DataRow 209 '11:46:46' '13:26' '15_ABASTO' 100 'LP1912'
Set-DataRow 218 '12:35:30' '13:50' '215A_EL PATO' 75 'LP1912'
Set-DataRow 219 '12:50:41' '13:50' '11_ETCHEVERRY' 60 'LP1912'
Set-DataRow 261 '14:20:49' '15:36' '23_HERNANDEZ' 76 'LP1912'
Set-DataRow 262 '15:36:13' '15:36' '10_OLMOS' 0 'LP1912'
Set-DataRow 302 '14:49:07' '16:43' '225_GOMEZ' 114 'LP1912'
Set-DataRow 304 '14:49:07' '16:43' '16_P MOR-SANTA ANA' 114 'LP1912'
Set-DataRow 314 '15:36:13' '17:05' '215A_EL PATO' 89 'LP1912'
Set-DataRow 315 '16:20:15' '17:05' '23_HERNANDEZ' 45 'LP1912'
Set-DataRow 316 '16:53:02' '17:05' '11_ETCHEVERRY' 12 'LP1912'
Set-DataRow 317 '16:34:19' '17:10' '10_OLMOS' 36 'LP1912'
Set-DataRow 318 '16:34:19' '17:16' '11_ETCHEVERRY' 42 'LP1912'
Set-DataRow 319 '15:59:02' '17:17' '11_ETCHEVERRY' 78 'LP1912'
Set-DataRow 320 '16:45:34' '17:20' '16_SANTA ANA' 35 'LP1912'
Set-DataRow 321 '15:36:13' '17:21' '26_HERNANDEZ' 105 'LP1912'
Set-DataRow 322 '16:20:15' '17:21' '16_SANTA ANA' 61 'LP1912'
Set-DataRow 323 '15:36:13' '17:24' '84_COLONIA URQUIZA-ESC 49' 108 'LP1912'
Set-DataRow 324 '16:34:19' '17:28' '14_ABASTO' 54 'LP1912'
Set-DataRow 325 '16:53:02' '17:29' '14_ABASTO' 36 'LP1912'
Set-DataRow 326 '16:34:19' '17:31' '15_ABASTO' 57 'LP1912'
Set-DataRow 327 '16:45:34' '17:34' '23_HERNANDEZ' 49 'LP1912'
Set-DataRow 328 '16:53:02' '17:35' '23_HERNANDEZ' 42 'LP1912'
Set-DataRow 329 '16:20:15' '17:36' '27_EL RETIRO' 76 'LP1912'
Set-DataRow 330 '15:59:02' '17:37' '27_EL RETIRO' 98 'LP1912'
Set-DataRow 331 '15:59:02' '17:38' '17_ROMERO' 99 'LP1912'
Set-DataRow 332 '16:45:34' '17:38' '27_EL RETIRO' 53 'LP1912'
Set-DataRow 333 '16:34:19' '17:39' '27_EL RETIRO' 65 'LP1912'
Set-DataRow 334 '16:45:34' '17:40' '16_SANTA ANA' 55 'LP1912'
Set-DataRow 335 '15:59:02' '17:40' '215B_EL PATO' 101 'LP1912'
Set-DataRow 336 '16:34:19' '17:41' '16_SANTA ANA' 67 'LP1912'
Set-DataRow 337 '16:45:34' '17:45' '15_ABASTO' 60 'LP1912'
Set-DataRow 338 '16:34:19' '17:50' '16_P MOR-167 Y 521' 76 'LP1912'
Set-DataRow 339 '15:59:02' '17:51' '16_P MOR-167 Y 521' 112 'LP1912'
Set-DataRow 340 '15:59:02' '17:52' '81_EL PELIGRO' 113 'LP1912'
Set-DataRow 341 '16:20:15' '18:04' '17_ROMERO' 104 'LP1912'
Set-DataRow 342 '16:34:19' '18:21' '26_HERNANDEZ' 107 'LP1912'
Set-DataRow 343 '16:53:02' '18:22' '26_HERNANDEZ' 89 'LP1912'
Set-DataRow 344 '16:34:19' '18:28' '215C_EL PATO' 114 'LP1912'
Set-DataRow 345 '16:34:19' '18:32' '11X44_ETCHEVERRY' 118 'LP1912'
Set-DataRow 346 '16:53:02' '18:45' '14_ABASTO' 112 'LP1912'
Set-DataRow 347 '16:53:02' '18:48' '14X44_ABASTO' 115 'LP1912'
